$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new data rows right before the former row 17 (pushes the old
# rows 17-52 down to 20-55, carrying their formatting/styles along).
$ws.Rows("17:19").Insert()

# --- Row 17 ---------------------------------------------------------------
$ws.Cells.Item(17, 1).Value = 1
$ws.Cells.Item(17, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(17, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(17, 4).Value = 44607
$ws.Cells.Item(17, 5).Value = 15
$ws.Cells.Item(17, 6).Value = "Fruta"
$ws.Cells.Item(17, 7).Value = 100103
$ws.Cells.Item(17, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(17, 9).Value = 100103006
$ws.Cells.Item(17, 10).Value = "Nectarín"
$ws.Cells.Item(17, 11).Value = "June Pearl"
$ws.Cells.Item(17, 12).Value = "Primera"
$ws.Cells.Item(17, 13).Value = 200
$ws.Cells.Item(17, 14).Value = 21000
$ws.Cells.Item(17, 15).Value = 22000
$ws.Cells.Item(17, 16).Value = 21500
$ws.Cells.Item(17, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(17, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(17, 19).Value = 1194
$ws.Cells.Item(17, 20).Value = 18

# --- Row 18 ---------------------------------------------------------------
$ws.Cells.Item(18, 1).Value = 1
$ws.Cells.Item(18, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(18, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(18, 4).Value = 44607
$ws.Cells.Item(18, 5).Value = 15
$ws.Cells.Item(18, 6).Value = "Fruta"
$ws.Cells.Item(18, 7).Value = 100103
$ws.Cells.Item(18, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(18, 9).Value = 100103006
$ws.Cells.Item(18, 10).Value = "Nectarín"
$ws.Cells.Item(18, 11).Value = "Ruby Diamond"
$ws.Cells.Item(18, 12).Value = "Primera"
$ws.Cells.Item(18, 13).Value = 250
$ws.Cells.Item(18, 14).Value = 21000
$ws.Cells.Item(18, 15).Value = 22000
$ws.Cells.Item(18, 16).Value = 21500
$ws.Cells.Item(18, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(18, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(18, 19).Value = 1194
$ws.Cells.Item(18, 20).Value = 18

# --- Row 19 ---------------------------------------------------------------
$ws.Cells.Item(19, 1).Value = 1
$ws.Cells.Item(19, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(19, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(19, 4).Value = 44607
$ws.Cells.Item(19, 5).Value = 15
$ws.Cells.Item(19, 6).Value = "Fruta"
$ws.Cells.Item(19, 7).Value = 100103
$ws.Cells.Item(19, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(19, 9).Value = 100103006
$ws.Cells.Item(19, 10).Value = "Nectarín"
$ws.Cells.Item(19, 11).Value = "Venus"
$ws.Cells.Item(19, 12).Value = "Primera"
$ws.Cells.Item(19, 13).Value = 270
$ws.Cells.Item(19, 14).Value = 21000
$ws.Cells.Item(19, 15).Value = 22000
$ws.Cells.Item(19, 16).Value = 21500
$ws.Cells.Item(19, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(19, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(19, 19).Value = 1194
$ws.Cells.Item(19, 20).Value = 18
